$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.072990800604276
$ws.Range("D2").Value = 1.071399269052526
$ws.Range("E2").Value = 1.076711303264512
$ws.Range("F2").Value = 1.081107597432229
$ws.Range("I2").Value = 1.044980875815597
$ws.Range("J2").Value = 1.07790756131058
$ws.Range("K2").Value = 1.074096545804653
$ws.Range("L2").Value = 1.079394515767882
$ws.Range("M2").Value = 1.083779283101807
$ws.Range("N2").Value = 1.079438313685367

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.075361588182843
$ws.Range("D3").Value = 1.073296318464397
$ws.Range("E3").Value = 1.079036358024944
$ws.Range("F3").Value = 1.08326712903639
$ws.Range("I3").Value = 1.045601643070414
$ws.Range("J3").Value = 1.079930243247953
$ws.Range("K3").Value = 1.075807474884754
$ws.Range("L3").Value = 1.081533447217192
$ws.Range("M3").Value = 1.085753952377518
$ws.Range("N3").Value = 1.081463868063096

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.076888202086994
$ws.Range("D4").Value = 1.074516910576057
$ws.Range("E4").Value = 1.080533768585569
$ws.Range("F4").Value = 1.084657597418943
$ws.Range("I4").Value = 1.045998642726188
$ws.Range("J4").Value = 1.081231545673189
$ws.Range("K4").Value = 1.076907219679291
$ws.Range("L4").Value = 1.082910094256724
$ws.Range("M4").Value = 1.087024435152986
$ws.Range("N4").Value = 1.082767018487037

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.077528251587805
$ws.Range("D5").Value = 1.075028423420171
$ws.Range("E5").Value = 1.081161632173218
$ws.Range("F5").Value = 1.085240538021227
$ws.Range("I5").Value = 1.046164434611527
$ws.Range("J5").Value = 1.081776853242504
$ws.Range("K5").Value = 1.077367827269569
$ws.Range("L5").Value = 1.083487107781163
$ws.Range("M5").Value = 1.087556843897924
$ws.Range("N5").Value = 1.083313100455648

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.077635617985668
$ws.Range("D6").Value = 1.075114214419807
$ws.Range("E6").Value = 1.081266957793723
$ws.Range("F6").Value = 1.085338322859804
$ws.Range("I6").Value = 1.046192207278126
$ws.Range("J6").Value = 1.081868310658066
$ws.Range("K6").Value = 1.077445065079253
$ws.Range("L6").Value = 1.083583890632942
$ws.Range("M6").Value = 1.087646138886101
$ws.Range("N6").Value = 1.083404687751231

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.076896761235066
$ws.Range("D7").Value = 1.074523751771798
$ws.Range("E7").Value = 1.080542164549169
$ws.Range("F7").Value = 1.084665392984372
$ws.Range("I7").Value = 1.046000862376594
$ws.Range("J7").Value = 1.081238838964982
$ws.Range("K7").Value = 1.076913381075825
$ws.Range("L7").Value = 1.082917811080872
$ws.Range("M7").Value = 1.087031555865178
$ws.Range("N7").Value = 1.082774322136141

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.073793586920637
$ws.Range("D8").Value = 1.072041840392684
$ws.Range("E8").Value = 1.077498552283503
$ws.Range("F8").Value = 1.081838871464502
$ws.Range("I8").Value = 1.045191643702934
$ws.Range("J8").Value = 1.0785927134616
$ws.Range("K8").Value = 1.074676301817867
$ws.Range("L8").Value = 1.08011893036414
$ws.Range("M8").Value = 1.084448156906299
$ws.Range("N8").Value = 1.08012443883103

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06826642283627
$ws.Range("D9").Value = 1.067613810180741
$ws.Range("E9").Value = 1.072079409390682
$ws.Range("F9").Value = 1.076803666483231
$ws.Range("I9").Value = 1.043729265929159
$ws.Range("J9").Value = 1.073870697209203
$ws.Range("K9").Value = 1.070676594918023
$ws.Range("L9").Value = 1.075128615463373
$ws.Range("M9").Value = 1.079838640986246
$ws.Range("N9").Value = 1.075395716773931

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.064539223596386
$ws.Range("D10").Value = 1.064622908132459
$ws.Range("E10").Value = 1.068426408264737
$ws.Range("F10").Value = 1.073407778217123
$ws.Range("I10").Value = 1.042729001825527
$ws.Range("J10").Value = 1.070680472769378
$ws.Range("K10").Value = 1.067969263706533
$ws.Range("L10").Value = 1.071760015793497
$ws.Range("M10").Value = 1.076724839803492
$ws.Range("N10").Value = 1.072200961849477

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.062914617225178
$ws.Range("D11").Value = 1.06331809943095
$ws.Range("E11").Value = 1.066834476261716
$ws.Range("F11").Value = 1.071927501215167
$ws.Range("I11").Value = 1.042289672139772
$ws.Range("J11").Value = 1.069288514195387
$ws.Range("K11").Value = 1.066786799368026
$ws.Range("L11").Value = 1.070290913043939
$ws.Range("M11").Value = 1.075366333353396
$ws.Range("N11").Value = 1.070807026534651

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.062309503403344
$ws.Range("D12").Value = 1.062831930881007
$ws.Range("E12").Value = 1.066241583611099
$ws.Range("F12").Value = 1.071376135002996
$ws.Range("I12").Value = 1.042125535912837
$ws.Range("J12").Value = 1.068769843392243
$ws.Range("K12").Value = 1.066346010399001
$ws.Range("L12").Value = 1.069743599591129
$ws.Range("M12").Value = 1.07486014424925
$ws.Range("N12").Value = 1.070287619159474

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.062439378499275
$ws.Range("D13").Value = 1.062936284466663
$ws.Range("E13").Value = 1.066368833363196
$ws.Range("F13").Value = 1.071494474735243
$ws.Range("I13").Value = 1.042160786937347
$ws.Range("J13").Value = 1.068881174817442
$ws.Range("K13").Value = 1.066440632798001
$ws.Range("L13").Value = 1.069861074444047
$ws.Range("M13").Value = 1.074968795737603
$ws.Range("N13").Value = 1.070399108688067

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.062864632565158
$ws.Range("D14").Value = 1.06327794352988
$ws.Range("E14").Value = 1.066785500070503
$ws.Range("F14").Value = 1.07188195651056
$ws.Range("I14").Value = 1.042276124058675
$ws.Range("J14").Value = 1.069245674344333
$ws.Range("K14").Value = 1.066750395862071
$ws.Range("L14").Value = 1.070245705346236
$ws.Range("M14").Value = 1.075324524087795
$ws.Range("N14").Value = 1.070764125846095

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.063126423661596
$ws.Range("D15").Value = 1.063488250295072
$ws.Range("E15").Value = 1.067042011476146
$ws.Range("F15").Value = 1.072120493240131
$ws.Range("I15").Value = 1.042347060722473
$ws.Range("J15").Value = 1.069470036392558
$ws.Range("K15").Value = 1.066941041957613
$ws.Range("L15").Value = 1.070482472548993
$ws.Range("M15").Value = 1.075543489416227
$ws.Range("N15").Value = 1.070988806514163

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.064646813216611
$ws.Range("D16").Value = 1.064709295379057
$ws.Range("E16").Value = 1.068531841105037
$ws.Range("F16").Value = 1.07350580805432
$ws.Range("I16").Value = 1.042758026450537
$ws.Range("J16").Value = 1.070772625752012
$ws.Range("K16").Value = 1.068047522190784
$ws.Range("L16").Value = 1.071857290388093
$ws.Range("M16").Value = 1.076814780436707
$ws.Range("N16").Value = 1.072293245699917

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.06559760942154
$ws.Range("D17").Value = 1.065472590108013
$ws.Range("E17").Value = 1.069463615767462
$ws.Range("F17").Value = 1.074372112871243
$ws.Range("I17").Value = 1.04301414036915
$ws.Range("J17").Value = 1.071586842544691
$ws.Range("K17").Value = 1.068738836096489
$ws.Range("L17").Value = 1.072716837885929
$ws.Range("M17").Value = 1.077609462449628
$ws.Range("N17").Value = 1.073108618773834

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.066151163663962
$ws.Range("D18").Value = 1.065916871304312
$ws.Range("E18").Value = 1.070006127590938
$ws.Range("F18").Value = 1.074876468160185
$ws.Range("I18").Value = 1.043162929386635
$ws.Range("J18").Value = 1.072060744784632
$ws.Range("K18").Value = 1.069141089144235
$ws.Range("L18").Value = 1.073217189909245
$ws.Range("M18").Value = 1.078072005017721
$ws.Range("N18").Value = 1.073583194009312

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.06633973843435
$ws.Range("D19").Value = 1.066068202285297
$ws.Range("E19").Value = 1.07019094601872
$ws.Range("F19").Value = 1.075048281313979
$ws.Range("I19").Value = 1.043213561686759
$ws.Range("J19").Value = 1.072222162180441
$ws.Range("K19").Value = 1.069278082264397
$ws.Range("L19").Value = 1.073387627606987
$ws.Range("M19").Value = 1.078229554875127
$ws.Range("N19").Value = 1.073744840636338

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.065495704822117
$ws.Range("D20").Value = 1.065390792899751
$ws.Range("E20").Value = 1.069363746549342
$ws.Range("F20").Value = 1.074279264719705
$ws.Range("I20").Value = 1.042986723733701
$ws.Range("J20").Value = 1.071499590284023
$ws.Range("K20").Value = 1.068664766177037
$ws.Range("L20").Value = 1.072624721136899
$ws.Range("M20").Value = 1.077524302489527
$ws.Range("N20").Value = 1.073021242604947

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.062739452182637
$ws.Range("D21").Value = 1.063177375331917
$ws.Range("E21").Value = 1.066662846100039
$ws.Range("F21").Value = 1.071767895283429
$ws.Range("I21").Value = 1.042242186502152
$ws.Range("J21").Value = 1.06913838384875
$ws.Range("K21").Value = 1.066659222006103
$ws.Range("L21").Value = 1.07013248634302
$ws.Range("M21").Value = 1.075219814920883
$ws.Range("N21").Value = 1.070656682985702

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.060996839241147
$ws.Range("D22").Value = 1.061776985390642
$ws.Range("E22").Value = 1.064955525144747
$ws.Range("F22").Value = 1.070180047625395
$ws.Range("I22").Value = 1.041768563598322
$ws.Range("J22").Value = 1.067644313523762
$ws.Range("K22").Value = 1.065389157844942
$ws.Range("L22").Value = 1.068556102506975
$ws.Range("M22").Value = 1.07376173151812
$ws.Range("N22").Value = 1.069160490909481

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.061921564271205
$ws.Range("D23").Value = 1.062520200320308
$ws.Range("E23").Value = 1.065861493973108
$ws.Range("F23").Value = 1.071022650812149
$ws.Range("I23").Value = 1.042020167371517
$ws.Range("J23").Value = 1.068437263701417
$ws.Range("K23").Value = 1.066063319209307
$ws.Range("L23").Value = 1.069392682721067
$ws.Range("M23").Value = 1.074535572829546
$ws.Range("N23").Value = 1.069954567167341

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.065541754301559
$ws.Range("D24").Value = 1.065427756426643
$ws.Range("E24").Value = 1.069408876165425
$ws.Range("F24").Value = 1.0743212217227
$ws.Range("I24").Value = 1.042999113977875
$ws.Range("J24").Value = 1.071539018964133
$ws.Range("K24").Value = 1.068698238207524
$ws.Range("L24").Value = 1.072666347841553
$ws.Range("M24").Value = 1.077562785642164
$ws.Range("N24").Value = 1.073060727278305

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.069702598179955
$ws.Range("D25").Value = 1.068765252826035
$ws.Range("E25").Value = 1.07348728289371
$ws.Range("F25").Value = 1.078112095426578
$ws.Range("I25").Value = 1.044111729525686
$ws.Range("J25").Value = 1.075098712728854
$ws.Range("K25").Value = 1.071717662368455
$ws.Range("L25").Value = 1.076425898855231
$ws.Range("M25").Value = 1.081037329083564
$ws.Range("N25").Value = 1.076625476216475
